$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 71,3
$data[0,0] = 52.07346600331662
$data[0,1] = 34.32666666666668
$data[0,2] = 70
$data[1,0] = 0.2768159203980091
$data[1,1] = 0
$data[1,2] = 0.8133333333333337
$data[2,0] = 1.081956882255389
$data[2,1] = 0.2000000000000001
$data[2,2] = 2.206666666666668
$data[3,0] = 0.3865339966832492
$data[3,1] = 0
$data[3,2] = 1
$data[4,0] = 6.881724709784403
$data[4,1] = 3.153333333333334
$data[4,2] = 12.09333333333333
$data[5,0] = 46.50388059701481
$data[5,1] = 30.18666666666666
$data[5,2] = 63.59999999999999
$data[6,0] = 92.90633499170801
$data[6,1] = 85.07999999999998
$data[6,2] = 97.04000000000001
$data[7,0] = 5.713366500829183
$data[7,1] = 2.413333333333333
$data[7,2] = 9.880000000000001
$data[8,0] = 92.7742288557213
$data[8,1] = 85.06666666666666
$data[8,2] = 96.93333333333334
$data[9,0] = 75.61917081260357
$data[9,1] = 61.37333333333331
$data[9,2] = 87.20666666666664
$data[10,0] = 27.61071310116076
$data[10,1] = 15.83333333333334
$data[10,2] = 40.86666666666667
$data[11,0] = 90.51721393034818
$data[11,1] = 82.30666666666664
$data[11,2] = 95.63333333333331
$data[12,0] = 85.78766169154218
$data[12,1] = 75.11999999999998
$data[12,2] = 93.04000000000001
$data[13,0] = 0.4273300165837469
$data[13,1] = 0
$data[13,2] = 1.04
$data[14,0] = 7.007927031509111
$data[14,1] = 3.266666666666666
$data[14,2] = 12.32666666666667
$data[15,0] = 27.33140961857368
$data[15,1] = 15.98
$data[15,2] = 39.52
$data[16,0] = 1.081956882255389
$data[16,1] = 0.2000000000000001
$data[16,2] = 2.206666666666668
$data[17,0] = 77.18195688225528
$data[17,1] = 63.90666666666667
$data[17,2] = 87.69333333333333
$data[18,0] = 8.376351575456042
$data[18,1] = 4.386666666666665
$data[18,2] = 13.56
$data[19,0] = 1.732802653399668
$data[19,1] = 0.5933333333333333
$data[19,2] = 3.346666666666668
$data[20,0] = 92.88470978441119
$data[20,1] = 84.57999999999998
$data[20,2] = 97.11333333333332
$data[21,0] = 0.4887893864013258
$data[21,1] = 0
$data[21,2] = 1.253333333333333
$data[22,0] = 6.979800995024868
$data[22,1] = 3.266666666666666
$data[22,2] = 12.23333333333334
$data[23,0] = 5.03389718076284
$data[23,1] = 2.52
$data[23,2] = 8.506666666666662
$data[24,0] = 1.758142620232171
$data[24,1] = 0.5866666666666666
$data[24,2] = 3.213333333333335
$data[25,0] = 90.96185737976771
$data[25,1] = 82.38000000000002
$data[25,2] = 96.12666666666669
$data[26,0] = 50.42782752902151
$data[26,1] = 32.87999999999999
$data[26,2] = 66.59333333333331
$data[27,0] = 0.5404311774461018
$data[27,1] = 0
$data[27,2] = 1.286666666666667
$data[28,0] = 25.19943615257042
$data[28,1] = 14.1
$data[28,2] = 38.48666666666666
$data[29,0] = 85.17276948590373
$data[29,1] = 74.26666666666667
$data[29,2] = 92.68000000000002
$data[30,0] = 90.68338308457699
$data[30,1] = 82.99999999999999
$data[30,2] = 95.4466666666667
$data[31,0] = 85.77615257048083
$data[31,1] = 74.41333333333331
$data[31,2] = 92.87333333333336
$data[32,0] = 58.0420895522387
$data[32,1] = 40.3
$data[32,2] = 73.95999999999998
$data[33,0] = 7.42666666666666
$data[33,1] = 3.639999999999999
$data[33,2] = 12.48
$data[34,0] = 46.52215588723045
$data[34,1] = 30.01333333333334
$data[34,2] = 62.94666666666664
$data[35,0] = 2.010082918739633
$data[35,1] = 0.6066666666666667
$data[35,2] = 3.673333333333332
$data[36,0] = 86.18348258706462
$data[36,1] = 76.95333333333332
$data[36,2] = 93.45333333333338
$data[37,0] = 28.27900497512423
$data[37,1] = 16.32666666666668
$data[37,2] = 42.70666666666666
$data[38,0] = 1.360398009950248
$data[38,1] = 0.2733333333333334
$data[38,2] = 2.7
$data[39,0] = 42.78245439469308
$data[39,1] = 26.50000000000001
$data[39,2] = 59.10000000000001
$data[40,0] = 76.54262023217237
$data[40,1] = 63.05333333333336
$data[40,2] = 87.48666666666666
$data[41,0] = 25.23734660033162
$data[41,1] = 14.13333333333334
$data[41,2] = 38.48666666666666
$data[42,0] = 70.13492537313427
$data[42,1] = 54.20666666666668
$data[42,2] = 82.15333333333332
$data[43,0] = 77.37167495854052
$data[43,1] = 65.31333333333335
$data[43,2] = 87.08666666666667
$data[44,0] = 84.98278606965157
$data[44,1] = 74.27333333333333
$data[44,2] = 91.84666666666668
$data[45,0] = 92.9890547263681
$data[45,1] = 84.33999999999997
$data[45,2] = 97.11999999999999
$data[46,0] = 27.81018242122709
$data[46,1] = 15.89333333333334
$data[46,2] = 40.86
$data[47,0] = 1.210049751243781
$data[47,1] = 0.4000000000000002
$data[47,2] = 2.320000000000001
$data[48,0] = 2.059336650082917
$data[48,1] = 0.6666666666666667
$data[48,2] = 3.859999999999999
$data[49,0] = 0.2844776119402977
$data[49,1] = 0
$data[49,2] = 0.8333333333333337
$data[50,0] = 84.22298507462685
$data[50,1] = 73.69999999999999
$data[50,2] = 91.67999999999999
$data[51,0] = 27.90723051409605
$data[51,1] = 16.17333333333334
$data[51,2] = 41.00666666666667
$data[52,0] = 0.42898839137645
$data[52,1] = 0
$data[52,2] = 1.04
$data[53,0] = 1.73927031509121
$data[53,1] = 0.5933333333333333
$data[53,2] = 3.246666666666668
$data[54,0] = 85.68315091210609
$data[54,1] = 75.10666666666664
$data[54,2] = 92.73333333333335
$data[55,0] = 90.33379767827522
$data[55,1] = 81.64000000000001
$data[55,2] = 95.70666666666666
$data[56,0] = 1.999237147595354
$data[56,1] = 0.6
$data[56,2] = 3.666666666666665
$data[57,0] = 73.75804311774455
$data[57,1] = 57.73999999999997
$data[57,2] = 85.7
$data[58,0] = 0.5696185737976766
$data[58,1] = 0
$data[58,2] = 1.373333333333333
$data[59,0] = 86.10195688225529
$data[59,1] = 75.05333333333333
$data[59,2] = 93.26666666666668
$data[60,0] = 75.58560530679925
$data[60,1] = 61.11333333333331
$data[60,2] = 87.16666666666664
$data[61,0] = 86.27880597014912
$data[61,1] = 74.7
$data[61,2] = 93.55333333333328
$data[62,0] = 0.5758872305140947
$data[62,1] = 0
$data[62,2] = 1.373333333333333
$data[63,0] = 89.93485903814243
$data[63,1] = 80.58666666666669
$data[63,2] = 95.11333333333337
$data[64,0] = 56.87870646766162
$data[64,1] = 39.69333333333331
$data[64,2] = 72.72666666666665
$data[65,0] = 54.75339966832492
$data[65,1] = 37.35333333333333
$data[65,2] = 71.11333333333337
$data[66,0] = 1.579800995024875
$data[66,1] = 0.4133333333333334
$data[66,2] = 3.046666666666667
$data[67,0] = 93.08036484245432
$data[67,1] = 85.90000000000001
$data[67,2] = 97.08666666666666
$data[68,0] = 1.376517412935322
$data[68,1] = 0.3400000000000001
$data[68,2] = 2.746666666666666
$data[69,0] = 16.67555555555554
$data[69,1] = 8.506666666666666
$data[69,2] = 26.06666666666666
$data[70,0] = 20.76772802653386
$data[70,1] = 10.74666666666667
$data[70,2] = 32.77333333333333

$range = $ws.Range("A2:C72")
$range.Value2 = $data
